$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.821893095970154
$ws.Range("B1").Value = 2.100561618804932
$ws.Range("C1").Value = 1.954269766807556
$ws.Range("D1").Value = 1.201654434204102
$ws.Range("E1").Value = 0.7704174518585205
